$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Row 4: End Year 2050 -> 2030
$ws.Range("B4").Value = 2030

# Row 16: realistic_candidate_capacities -> realistic_candidate_capacities_tobe_installed, FALSE -> TRUE, description update
$ws.Range("A16").Value = "realistic_candidate_capacities_tobe_installed"
$ws.Range("B16").Value = $true
$ws.Range("C16").Value = "If this is true, the real capacity  of the power plants is chosen"

# Row 17: description update (name/value unchanged)
$ws.Range("C17").Value = "If this is true, the real capacity of the candidate power plants is considered for the future investments. Otherwise the dummy capacity inidated "

# Column A width update (target stored width 48.54296875; engine quantizes ColumnWidth
# to a 1/6-character pixel grid, so 47.666666666666664 is the closest achievable input)
$ws.Columns("A").ColumnWidth = 47.666666666666664

# Selection update
$ws.Range("B5").Select()
